$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 165: base_image_url for "guarda la fecha" (underscore -> hyphen filename)
$ws.Range("G165").Value = "/media/product_images/invitaciones_papeleria/bodas/guarda_la_fecha/guarda-la-fecha.jpg"

# --- New row 167: Carteles de carton espuma
$ws.Range("A167").Value = "carteles-carton-espuma"
$ws.Range("B167").Value = "Carteles de cartón espuma"
$ws.Range("C167").Value = "invitaciones-papeleria"
$ws.Range("D167").Value = "bodas"
$ws.Range("E167").Value = "CCE-001"
$ws.Range("F167").Value = "Carteles para tu día especial"
$ws.Range("G167").Value = "/media/product_images/invitaciones_papeleria/bodas/carteles_carton_espuma/carteles-carton-espuma.jpg"
$ws.Range("H167").Value = "active"
$ws.Range("I167").Value = "unisex"

# --- New row 168: Invitaciones de boda
$ws.Range("A168").Value = "invitaciones-boda"
$ws.Range("B168").Value = "Invitaciones para tu boda"
$ws.Range("C168").Value = "invitaciones-papeleria"
$ws.Range("D168").Value = "bodas"
$ws.Range("E168").Value = "IB-001"
$ws.Range("F168").Value = "Invitaciones elegantes para tu boda"
$ws.Range("G168").Value = "/media/product_images/invitaciones_papeleria/bodas/invitaciones_de_boda/invitaciones-de-boda.jpg"
$ws.Range("H168").Value = "active"
$ws.Range("I168").Value = "unisex"

# --- New row 169: Invitaciones despedida de soltera
$ws.Range("A169").Value = "invitaciones-despedida-de-soltera"
$ws.Range("B169").Value = "Invitaciones para tu despedida de soltera"
$ws.Range("C169").Value = "invitaciones-papeleria"
$ws.Range("D169").Value = "bodas"
$ws.Range("E169").Value = "IDS-001"
$ws.Range("F169").Value = "Invitaciones para tu despedida de soltera"
$ws.Range("G169").Value = "/media/product_images/invitaciones_papeleria/bodas/invitaciones_despedida_de_soltera/invitaciones-despedida-de-soltera.jpg"
$ws.Range("H169").Value = "active"
$ws.Range("I169").Value = "unisex"

# --- New row 170: Programa de boda
$ws.Range("A170").Value = "programa-de-boda"
$ws.Range("B170").Value = "Programa de boda"
$ws.Range("C170").Value = "invitaciones-papeleria"
$ws.Range("D170").Value = "bodas"
$ws.Range("E170").Value = "PDB-001"
$ws.Range("F170").Value = "Programa de boda"
$ws.Range("G170").Value = "/media/product_images/invitaciones_papeleria/bodas/programa_de_boda/programa-de-boda.jpg"
$ws.Range("H170").Value = "active"
$ws.Range("I170").Value = "unisex"

# --- New row 171: Libro de firmas
$ws.Range("A171").Value = "libro-de-firmas"
$ws.Range("B171").Value = "Libro de firmas de invitados"
$ws.Range("C171").Value = "invitaciones-papeleria"
$ws.Range("D171").Value = "bodas"
$ws.Range("E171").Value = "LDF-001"
$ws.Range("F171").Value = "Libro de firmas de invitados"
$ws.Range("G171").Value = "/media/product_images/invitaciones_papeleria/bodas/libro_de_firmas/libro-de-firmas.jpg"
$ws.Range("H171").Value = "active"
$ws.Range("I171").Value = "unisex"

# --- New row 172: Tarjetas de menu
$ws.Range("A172").Value = "tarjetas-de-menu"
$ws.Range("B172").Value = "Tarjetas de menú"
$ws.Range("C172").Value = "invitaciones-papeleria"
$ws.Range("D172").Value = "bodas"
$ws.Range("E172").Value = "TM-001"
$ws.Range("F172").Value = "Tarjetas para el menú de la boda"
$ws.Range("G172").Value = "/media/product_images/invitaciones_papeleria/bodas/tarjetas_de_menu/tarjetas-de-menu.jpg"
$ws.Range("H172").Value = "active"
$ws.Range("I172").Value = "unisex"

# --- Column G width (approximates Excel's auto best-fit width after the longer URLs were added)
$ws.Columns.Item(7).ColumnWidth = 129.83

# --- Restore view/selection state as close as possible to the saved workbook
$excel.ActiveWindow.ScrollRow = 149
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G169").Select()
